$d = $word.ActiveDocument

# Remove the trailing "Ver no Jupiter ..." and "(c) 2020 ..." paragraphs,
# plus the now-redundant blank paragraph that used to separate them from
# the bibliography entry above, leaving the single blank paragraph that
# precedes the final page-break paragraph untouched.
$startPara = $d.Paragraphs.Item(38)
$endPara   = $d.Paragraphs.Item(40)

$range = $d.Range($startPara.Range.Start, $endPara.Range.End)
$range.Delete()
